$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 - this shifts rows 5:9 down to 6:10
$ws.Range("A5:R5").Insert()

# Copy style (date format) from the row above insertion point into new D5 - use same style as other D cells
$ws.Range("D6").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats

# Fill new row 5 with the same static data as the rest of the rows, differing only in D, J, K, L, M, P
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(5, 3).Value = "Los Lagos"
$ws.Cells.Item(5, 4).Value = 44526
$ws.Cells.Item(5, 5).Value = 10
$ws.Cells.Item(5, 6).Value = 100112012
$ws.Cells.Item(5, 7).Value = "Espinaca"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 25
$ws.Cells.Item(5, 11).Value = 9000
$ws.Cells.Item(5, 12).Value = 9000
$ws.Cells.Item(5, 13).Value = 9000
$ws.Cells.Item(5, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(5, 15).Value = "Región Metropolitana"
$ws.Cells.Item(5, 16).Value = 900
$ws.Cells.Item(5, 17).Value = 10
$ws.Cells.Item(5, 18).Value = "Hortaliza"
